$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 1325.0834
$ws.Range("I62").Value = 1381
$ws.Range("J62").Value = 1285.1428
$ws.Range("K62").Value = 1381
$ws.Range("L62").Value = 1285.1428
$ws.Range("M62").Value = -757
$ws.Range("N62").Value = -2533.1428
# Row 64
$ws.Range("H64").Value = 3961.2903
$ws.Range("I64").Value = 3895.8333
$ws.Range("J64").Value = 4185.7144
$ws.Range("K64").Value = 3895.8333
$ws.Range("L64").Value = 4185.7144
$ws.Range("M64").Value = -3647.8333
$ws.Range("N64").Value = -4681.7144
# Row 65
$ws.Range("H65").Value = 1325.0834
$ws.Range("I65").Value = 1381
$ws.Range("J65").Value = 1285.1428
$ws.Range("K65").Value = 6905
$ws.Range("L65").Value = 6425.714
$ws.Range("M65").Value = -3785
$ws.Range("N65").Value = -12665.714
# Row 67
$ws.Range("H67").Value = 3961.2903
$ws.Range("I67").Value = 3895.8333
$ws.Range("J67").Value = 4185.7144
$ws.Range("K67").Value = 3895.8333
$ws.Range("L67").Value = 4185.7144
$ws.Range("M67").Value = -3037.8333
$ws.Range("N67").Value = -5901.7144
# Row 87
$ws.Range("H87").Value = 29900
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 29900
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 29900
$ws.Range("N87").Value = -32396
# Row 90
$ws.Range("H90").Value = 29900
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 29900
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 89700
$ws.Range("N90").Value = -102180
# Row 132
$ws.Range("H132").Value = 574.1803
$ws.Range("I132").Value = 517.9825
$ws.Range("J132").Value = 1375
$ws.Range("K132").Value = 1553.9475
$ws.Range("L132").Value = 4125
$ws.Range("M132").Value = 976.0525000000002
$ws.Range("N132").Value = -9185
# Row 137
$ws.Range("H137").Value = 949.1385
$ws.Range("I137").Value = 823.13336
$ws.Range("J137").Value = 2461.2
$ws.Range("K137").Value = 2469.40008
$ws.Range("L137").Value = 7383.599999999999
$ws.Range("M137").Value = 80.59991999999966
$ws.Range("N137").Value = -12483.6

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1851.1364
$ws.Range("I2").Value = 1371.5294
$ws.Range("J2").Value = 3481.8
$ws.Range("K2").Value = 1371.5294
$ws.Range("L2").Value = 3481.8
$ws.Range("M2").Value = -1258.5294
$ws.Range("N2").Value = -3707.8
# Row 5
$ws.Range("H5").Value = 349.75
$ws.Range("I5").Value = 349.75
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 349.75
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -237.75
# Row 32
$ws.Range("H32").Value = 4967
$ws.Range("I32").Value = 3662.6875
$ws.Range("J32").Value = 10658.546
$ws.Range("K32").Value = 3662.6875
$ws.Range("L32").Value = 10658.546
$ws.Range("M32").Value = -3375.6875
$ws.Range("N32").Value = -11232.546
# Row 45
$ws.Range("H45").Value = 5043.4614
$ws.Range("I45").Value = 5175.2173
$ws.Range("J45").Value = 4033.3333
$ws.Range("K45").Value = 5175.2173
$ws.Range("L45").Value = 4033.3333
$ws.Range("M45").Value = -4798.2173
$ws.Range("N45").Value = -4787.3333
# Row 61
$ws.Range("H61").Value = 4894.1
$ws.Range("I61").Value = 5464.44
$ws.Range("J61").Value = 2042.4
$ws.Range("K61").Value = 5464.44
$ws.Range("L61").Value = 2042.4
$ws.Range("M61").Value = -5252.44
$ws.Range("N61").Value = -2466.4
# Row 74
$ws.Range("H74").Value = 1412.4828
$ws.Range("I74").Value = 1262
$ws.Range("J74").Value = 1885.4286
$ws.Range("K74").Value = 1262
$ws.Range("L74").Value = 1885.4286
$ws.Range("M74").Value = -388
$ws.Range("N74").Value = -3633.4286
# Row 77
$ws.Range("H77").Value = 1412.4828
$ws.Range("I77").Value = 1262
$ws.Range("J77").Value = 1885.4286
$ws.Range("K77").Value = 6310
$ws.Range("L77").Value = 9427.143
$ws.Range("M77").Value = -1942
$ws.Range("N77").Value = -18163.143
# Row 88
$ws.Range("H88").Value = 2886.1428
$ws.Range("I88").Value = 2798.6667
$ws.Range("J88").Value = 2951.75
$ws.Range("K88").Value = 2798.6667
$ws.Range("L88").Value = 2951.75
$ws.Range("M88").Value = -2392.6667
$ws.Range("N88").Value = -3763.75
# Row 91
$ws.Range("H91").Value = 2886.1428
$ws.Range("I91").Value = 2798.6667
$ws.Range("J91").Value = 2951.75
$ws.Range("K91").Value = 2798.6667
$ws.Range("L91").Value = 2951.75
$ws.Range("M91").Value = -1394.6667
$ws.Range("N91").Value = -5759.75
# Row 116
$ws.Range("H116").Value = 1851.1364
$ws.Range("I116").Value = 1371.5294
$ws.Range("J116").Value = 3481.8
$ws.Range("K116").Value = 1371.5294
$ws.Range("L116").Value = 3481.8
$ws.Range("M116").Value = 922.4706000000001
$ws.Range("N116").Value = -8069.8
# Row 122
$ws.Range("H122").Value = 988634.25
$ws.Range("I122").Value = 1284788
$ws.Range("J122").Value = 1455
$ws.Range("K122").Value = 3854364
$ws.Range("L122").Value = 4365
$ws.Range("M122").Value = -3851914
$ws.Range("N122").Value = -9265
# Row 132
$ws.Range("H132").Value = 2814.3333
$ws.Range("I132").Value = 1789.8462
$ws.Range("J132").Value = 4216.263
$ws.Range("K132").Value = 5369.5386
$ws.Range("L132").Value = 12648.789
$ws.Range("M132").Value = -2839.5386
$ws.Range("N132").Value = -17708.789
# Row 136
$ws.Range("H136").Value = 4894.1
$ws.Range("I136").Value = 5464.44
$ws.Range("J136").Value = 2042.4
$ws.Range("K136").Value = 16393.32
$ws.Range("L136").Value = 6127.200000000001
$ws.Range("M136").Value = -13843.32
$ws.Range("N136").Value = -11227.2

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1851.1364
$ws.Range("I3").Value = 1371.5294
$ws.Range("J3").Value = 3481.8
$ws.Range("K3").Value = 1371.5294
$ws.Range("L3").Value = 3481.8
$ws.Range("M3").Value = -1257.5294
$ws.Range("N3").Value = -3709.8
# Row 4
$ws.Range("H4").Value = 349.75
$ws.Range("I4").Value = 349.75
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 349.75
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -234.75
# Row 86
$ws.Range("H86").Value = 8335036.5
$ws.Range("I86").Value = 14494385
$ws.Range("J86").Value = 1799.4706
$ws.Range("K86").Value = 14494385
$ws.Range("L86").Value = 1799.4706
$ws.Range("M86").Value = -14493262
$ws.Range("N86").Value = -4045.4706
# Row 89
$ws.Range("H89").Value = 8335036.5
$ws.Range("I89").Value = 14494385
$ws.Range("J89").Value = 1799.4706
$ws.Range("K89").Value = 72471925
$ws.Range("L89").Value = 8997.353000000001
$ws.Range("M89").Value = -72466309
$ws.Range("N89").Value = -20229.353
# Row 99
$ws.Range("H99").Value = 125001050
$ws.Range("I99").Value = 250000130
$ws.Range("J99").Value = 1967
$ws.Range("K99").Value = 250000130
$ws.Range("L99").Value = 1967
$ws.Range("M99").Value = -249998632
$ws.Range("N99").Value = -4963
# Row 134
$ws.Range("H134").Value = 5383.2583
$ws.Range("I134").Value = 6669.857
$ws.Range("J134").Value = 2681.4
$ws.Range("K134").Value = 20009.571
$ws.Range("L134").Value = 8044.200000000001
$ws.Range("M134").Value = -17474.571
$ws.Range("N134").Value = -13114.2

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2897.34
$ws.Range("I31").Value = 1379.6072
$ws.Range("J31").Value = 4829
$ws.Range("K31").Value = 1379.6072
$ws.Range("L31").Value = 4829
$ws.Range("M31").Value = -1084.6072
$ws.Range("N31").Value = -5419
# Row 34
$ws.Range("H34").Value = 2897.34
$ws.Range("I34").Value = 1379.6072
$ws.Range("J34").Value = 4829
$ws.Range("K34").Value = 1379.6072
$ws.Range("L34").Value = 4829
$ws.Range("M34").Value = -1177.6072
$ws.Range("N34").Value = -5233
# Row 58
$ws.Range("H58").Value = 1402.8889
$ws.Range("I58").Value = 1032.8572
$ws.Range("J58").Value = 1920.9333
$ws.Range("K58").Value = 1032.8572
$ws.Range("L58").Value = 1920.9333
$ws.Range("M58").Value = -829.8571999999999
$ws.Range("N58").Value = -2326.9333
# Row 132
$ws.Range("H132").Value = 2043.1333
$ws.Range("I132").Value = 1884.0625
$ws.Range("J132").Value = 2434.6924
$ws.Range("K132").Value = 5652.1875
$ws.Range("L132").Value = 7304.0772
$ws.Range("M132").Value = -3122.1875
$ws.Range("N132").Value = -12364.0772
# Row 134
$ws.Range("H134").Value = 2388.5454
$ws.Range("I134").Value = 2537.3333
$ws.Range("J134").Value = 1719
$ws.Range("K134").Value = 7611.999899999999
$ws.Range("L134").Value = 5157
$ws.Range("M134").Value = -5076.999899999999
$ws.Range("N134").Value = -10227
# Row 136
$ws.Range("H136").Value = 1402.8889
$ws.Range("I136").Value = 1032.8572
$ws.Range("J136").Value = 1920.9333
$ws.Range("K136").Value = 3098.5716
$ws.Range("L136").Value = 5762.7999
$ws.Range("M136").Value = -548.5715999999998
$ws.Range("N136").Value = -10862.7999

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1299639.1
$ws.Range("I131").Value = 5882956.5
$ws.Range("J131").Value = 1032.6333
$ws.Range("K131").Value = 17648869.5
$ws.Range("L131").Value = 3097.8999
$ws.Range("M131").Value = -17643829.5
$ws.Range("N131").Value = -13177.8999

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -888
$ws.Range("N5").ClearContents()
# Row 113
$ws.Range("H113").Value = 71429920
$ws.Range("I113").Value = 142858000
$ws.Range("J113").Value = 1835.7142
$ws.Range("K113").Value = 142858000
$ws.Range("L113").Value = 1835.7142
$ws.Range("M113").Value = -142855830
$ws.Range("N113").Value = -6175.7142
# Row 132
$ws.Range("H132").Value = 2591.4912
$ws.Range("I132").Value = 2558.9722
$ws.Range("J132").Value = 2647.238
$ws.Range("K132").Value = 7676.9166
$ws.Range("L132").Value = 7941.714
$ws.Range("M132").Value = -5146.9166
$ws.Range("N132").Value = -13001.714

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 657.7059
$ws.Range("I16").Value = 657.0714
$ws.Range("J16").Value = 660.6667
$ws.Range("K16").Value = 657.0714
$ws.Range("L16").Value = 660.6667
$ws.Range("M16").Value = -487.0714
$ws.Range("N16").Value = -1000.6667
# Row 46
$ws.Range("H46").Value = 19608820
$ws.Range("I46").Value = 30303814
$ws.Range("J46").Value = 1333.3334
$ws.Range("K46").Value = 30303814
$ws.Range("L46").Value = 1333.3334
$ws.Range("M46").Value = -30303626
$ws.Range("N46").Value = -1709.3334
# Row 132
$ws.Range("H132").Value = 12392633
$ws.Range("I132").Value = 19103830
$ws.Range("J132").Value = 2729.7693
$ws.Range("K132").Value = 57311490
$ws.Range("L132").Value = 8189.3079
$ws.Range("M132").Value = -57308960
$ws.Range("N132").Value = -13249.3079
# Row 136
$ws.Range("H136").Value = 9204.280000000001
$ws.Range("I136").Value = 13954.7
$ws.Range("J136").Value = 6037.3335
$ws.Range("K136").Value = 41864.10000000001
$ws.Range("L136").Value = 18112.0005
$ws.Range("M136").Value = -39314.10000000001
$ws.Range("N136").Value = -23212.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 60000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 60000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 60000
$ws.Range("N2").Value = -60224
$ws.Range("M2").ClearContents()
# Row 96
$ws.Range("H96").Value = 2368.6365
$ws.Range("I96").Value = 1751.6666
$ws.Range("J96").Value = 3109
$ws.Range("K96").Value = 1751.6666
$ws.Range("L96").Value = 3109
$ws.Range("M96").Value = -378.6666
$ws.Range("N96").Value = -5855
# Row 132
$ws.Range("H132").Value = 1612.1714
$ws.Range("I132").Value = 1005.63635
$ws.Range("J132").Value = 2638.6155
$ws.Range("K132").Value = 3016.90905
$ws.Range("L132").Value = 7915.8465
$ws.Range("M132").Value = -486.9090500000002
$ws.Range("N132").Value = -12975.8465
# Row 136
$ws.Range("H136").Value = 2671.8333
$ws.Range("I136").Value = 3316.3684
$ws.Range("J136").Value = 1951.4706
$ws.Range("K136").Value = 9949.1052
$ws.Range("L136").Value = 5854.4118
$ws.Range("M136").Value = -7399.1052
$ws.Range("N136").Value = -10954.4118
